$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 14-17 (sending cluster "Resolving-Mac" block) entirely.
$ws.Range("A14:T17").EntireRow.Delete()

# Columns B (Ligand symbol) and C (Receptor symbol) stay "Vtn"/"Itgav" for every
# remaining data row (2-13); only the underlying shared-string index changes,
# which Excel manages automatically, so no value change is actually needed here.
# Column D (Target cluster) for the 4th row of each sending-cluster block switches
# from "Resolving-Mac" to itself using the refreshed TPM numbers below.

$ws.Range("G2").Value = 7.844453333333334
$ws.Range("H2").Value = 23.53336
$ws.Range("I2").Value = 0.1489290605659587
$ws.Range("J2").Value = 0.1489290605659588
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 69.19315637612445
$ws.Range("R2").Value = 622.73840738512
$ws.Range("S2").Value = 0.009554017655279865
$ws.Range("T2").Value = 0.009554017655279867

$ws.Range("G3").Value = 7.844453333333334
$ws.Range("H3").Value = 23.53336
$ws.Range("I3").Value = 0.1489290605659587
$ws.Range("J3").Value = 0.1489290605659588
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 429.1820831025333
$ws.Range("R3").Value = 3862.6387479228
$ws.Range("S3").Value = 0.05926038663422309
$ws.Range("T3").Value = 0.0592603866342231

$ws.Range("G4").Value = 7.844453333333334
$ws.Range("H4").Value = 23.53336
$ws.Range("I4").Value = 0.1489290605659587
$ws.Range("J4").Value = 0.1489290605659588
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 171.8575648873778
$ws.Range("R4").Value = 1546.7180839864
$ws.Range("S4").Value = 0.02372966193653757
$ws.Range("T4").Value = 0.02372966193653758

$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 7.844453333333334
$ws.Range("H5").Value = 23.53336
$ws.Range("I5").Value = 0.1489290605659587
$ws.Range("J5").Value = 0.1489290605659588
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 408.3576011054133
$ws.Range("R5").Value = 3675.21840994872
$ws.Range("S5").Value = 0.05638499433991818
$ws.Range("T5").Value = 0.05638499433991821

$ws.Range("G6").Value = 20.35396833333334
$ws.Range("H6").Value = 61.06190500000001
$ws.Range("I6").Value = 0.3864255740794268
$ws.Range("J6").Value = 0.3864255740794268
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 179.5351765021678
$ws.Range("R6").Value = 1615.81658851951
$ws.Range("S6").Value = 0.02478976731053372
$ws.Range("T6").Value = 0.02478976731053372

$ws.Range("G7").Value = 20.35396833333334
$ws.Range("H7").Value = 61.06190500000001
$ws.Range("I7").Value = 0.3864255740794268
$ws.Range("J7").Value = 0.3864255740794268
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("Q7").Value = 1113.596850858059
$ws.Range("R7").Value = 10022.37165772253
$ws.Range("S7").Value = 0.1537626628293707
$ws.Range("T7").Value = 0.1537626628293708

$ws.Range("G8").Value = 20.35396833333334
$ws.Range("H8").Value = 61.06190500000001
$ws.Range("I8").Value = 0.3864255740794268
$ws.Range("J8").Value = 0.3864255740794268
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 445.9180627281612
$ws.Range("R8").Value = 4013.262564553451
$ws.Range("S8").Value = 0.0615712487656235
$ws.Range("T8").Value = 0.0615712487656235

$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 20.35396833333334
$ws.Range("H9").Value = 61.06190500000001
$ws.Range("I9").Value = 0.3864255740794268
$ws.Range("J9").Value = 0.3864255740794268
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 1059.563659618799
$ws.Range("R9").Value = 9536.072936569186
$ws.Range("S9").Value = 0.1463018951738988
$ws.Range("T9").Value = 0.1463018951738988

$ws.Range("G10").Value = 24.47399366666667
$ws.Range("H10").Value = 73.421981
$ws.Range("I10").Value = 0.4646453653546145
$ws.Range("J10").Value = 0.4646453653546145
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 215.8764669719002
$ws.Range("R10").Value = 1942.888202747102
$ws.Range("S10").Value = 0.02980768163830506
$ws.Range("T10").Value = 0.02980768163830506

$ws.Range("G11").Value = 24.47399366666667
$ws.Range("H11").Value = 73.421981
$ws.Range("I11").Value = 0.4646453653546145
$ws.Range("J11").Value = 0.4646453653546145
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 1339.009761083612
$ws.Range("R11").Value = 12051.0878497525
$ws.Range("S11").Value = 0.1848871126566959
$ws.Range("T11").Value = 0.1848871126566959

$ws.Range("G12").Value = 24.47399366666667
$ws.Range("H12").Value = 73.421981
$ws.Range("I12").Value = 0.4646453653546145
$ws.Range("J12").Value = 0.4646453653546145
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 536.180250668299
$ws.Range("R12").Value = 4825.622256014691
$ws.Range("S12").Value = 0.07403442550663759
$ws.Range("T12").Value = 0.07403442550663759

$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 24.47399366666667
$ws.Range("H13").Value = 73.421981
$ws.Range("I13").Value = 0.4646453653546145
$ws.Range("J13").Value = 0.4646453653546145
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 1274.03923747256
$ws.Range("R13").Value = 11466.35313725304
$ws.Range("S13").Value = 0.1759161455529759
$ws.Range("T13").Value = 0.1759161455529759

Write-Host "edits applied"
